# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# worksheets, which carry duplicate data in this workbook.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1854
    4  = 110
    7  = 1537
    9  = 598
    13 = 89
    16 = 134
    19 = 3599
    21 = 322
    23 = 134
    26 = 1362
    27 = 138
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
